$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.494.93'
$ws.Range("E2").Value = '  -0.54%  '
$ws.Range("D3").Value = '2.185.65'
$ws.Range("E3").Value = '  -1.99%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '''253.72'
$ws.Range("E5").Value = '  +3.84%  '
$ws.Range("D6").Value = '''0.611'
$ws.Range("E6").Value = '  -1.17%  '
$ws.Range("D7").Value = '''74.19'
$ws.Range("E7").Value = '  -1.37%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").Value = '''0.582'
$ws.Range("E9").Value = '  -3.90%  '
$ws.Range("D10").Value = '''40.12'
$ws.Range("E10").Value = '  -2.84%  '
$ws.Range("D11").Value = '''0.0914'
$ws.Range("E11").Value = '  -1.81%  '
$ws.Range("E12").Value = '  -0.03%  '
$ws.Range("D13").Value = '''6.76'
$ws.Range("E13").Value = '  -2.30%  '
$ws.Range("D14").Value = '2.511.66'
$ws.Range("E14").Value = '  -1.97%  '
$ws.Range("D15").Value = '''14.20'
$ws.Range("E15").Value = '  -3.21%  '
$ws.Range("D16").Value = '2.184.34'
$ws.Range("E16").Value = '  -2.11%  '
$ws.Range("D17").Value = '''0.772'
$ws.Range("E17").Value = '  -5.24%  '
$ws.Range("D18").Value = '42.398.36'
$ws.Range("E18").Value = '  -0.62%  '
$ws.Range("E19").Value = '  -2.38%  '
$ws.Range("D20").Value = '''70.95'
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("E21").Value = '  -1.40%  '
$ws.Range("D22").Value = '''227.07'
$ws.Range("E22").Value = '  -1.30%  '
$ws.Range("D23").Value = '''9.43'
$ws.Range("E23").Value = '  -6.43%  '
$ws.Range("D24").Value = '''2.12'
$ws.Range("E24").Value = '  -1.76%  '
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("D26").Value = '''10.46'
$ws.Range("E26").Value = '  -4.30%  '
$ws.Range("D27").Value = '''3.39'
$ws.Range("E27").Value = '  +1.97%  '
$ws.Range("D28").Value = '''2.17'
$ws.Range("E28").Value = '  -2.80%  '
$ws.Range("D29").Value = '''2.19'
$ws.Range("E29").Value = '  -0.08%  '
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").Value = '''37.67'
$ws.Range("E30").Value = '  +1.02%  '
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").Value = '''171.92'
$ws.Range("E31").Value = '  -1.16%  '
$ws.Range("D32").Value = '''20.09'
$ws.Range("E32").Value = '  -1.29%  '
$ws.Range("E33").Value = '  +3.94%  '
$ws.Range("D34").Value = '''5.15'
$ws.Range("E35").Value = '  -1.54%  '
$ws.Range("D36").Value = '''0.108'
$ws.Range("E36").Value = '  +0.30%  '
$ws.Range("E37").Value = '  +3.02%  '
$ws.Range("E38").Value = '  -4.24%  '
$ws.Range("D39").Value = '''12.09'
$ws.Range("E39").Value = '  -7.40%  '
$ws.Range("D40").Value = '''2.06'
$ws.Range("E40").Value = '  -3.52%  '
$ws.Range("B41").Value = 'NEARProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D41").Value = '''2.60'
$ws.Range("E41").Value = '  +13.05%  '
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").Value = '''0.196'
$ws.Range("E42").Value = '  -0.74%  '
$ws.Range("B43").Value = 'MultiversX'
$ws.Range("C43").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D43").Value = '''59.09'
$ws.Range("E43").Value = '  -1.96%  '
$ws.Range("D44").Value = '''5.16'
$ws.Range("E44").Value = '  -6.66%  '
$ws.Range("D45").Value = '''101.66'
$ws.Range("E45").Value = '  +1.07%  '
$ws.Range("D46").Value = '''0.0976'
$ws.Range("E46").Value = '  -1.68%  '
$ws.Range("D47").Value = '''0.463'
$ws.Range("E47").Value = '  +3.63%  '
$ws.Range("D48").Value = '''8.20'
$ws.Range("E48").Value = '  -4.56%  '
$ws.Range("E49").Value = '  -2.11%  '
$ws.Range("E50").Value = '  -1.66%  '
$ws.Range("E51").Value = '  +0.19%  '